## Add "Member assignment volume" table (columns J:L) mirroring the existing
## Name/上周/本周 layout in A:C, plus a second bar chart ("制造协同周任务量分析表")
## built from the new table - per commit message "add assignment volume function".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header relabel: "Period" -> "Name" (shared with the new table's header) ---
$ws.Range("A1").Value = "Name"

# --- New task/assignment-volume table in columns J:L, mirroring A:C's names ---
$ws.Range("J1").Value = "Name"
$ws.Range("K1").Value = "上周"
$ws.Range("L1").Value = "本周"

$ws.Range("J2").Value = "罗远明"
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 1

$ws.Range("J3").Value = "刘兴国"
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0

$ws.Range("J4").Value = "王超"
$ws.Range("K4").Value = 8
$ws.Range("L4").Value = 1

$ws.Range("J5").Value = "王言章"
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 3

$ws.Range("J6").Value = "李若愚"
$ws.Range("K6").Value = 5
$ws.Range("L6").Value = 4

$ws.Range("J7").Value = "聂正"
$ws.Range("K7").Value = 9
$ws.Range("L7").Value = 0

# --- New clustered-column chart sourced from J1:L7 ("制造协同周任务量分析表") ---
$co2 = $ws.ChartObjects().Add(0, 0, 300, 200)
$chart = $co2.Chart
$chart.ChartType = 51
$chart.SetSourceData($ws.Range("J1:L7"))

$chart.HasTitle = $true
$chart.ChartTitle.Text = "制造协同周任务量分析表"

$catAxis = $chart.Axes(1)
$catAxis.HasTitle = $true
$catAxis.AxisTitle.Text = "Member name"

$valAxis = $chart.Axes(2)
$valAxis.HasTitle = $true
$valAxis.AxisTitle.Text = "Weekly code volume data"

$chart.HasLegend = $true
$chart.Legend.Position = -4152

# Series border colors (上周 = yellow, 本周 = purple). Apply the 上周 series
# last so it wins on the shared styling slot.
$chart.SeriesCollection(2).Border.Color = 8388736
$chart.SeriesCollection(1).Border.Color = 65535

Write-Output "done"
